# AIP-68 AIP-273 Updated Test Data Added Fields for Voltage and current injection
#
# Adds two new columns of test data to Sheet1:
#   I: RMSInjectedVoltage = 60
#   J: RMSInjectedCurrent = 1
#
# The numeric-looking values "60" and "1" are entered the same way the
# existing PrefaultTime/PostFaultTime columns (G/H) were — as text via a
# leading apostrophe — so they round-trip as shared-string text cells with
# the workbook's existing "quote prefixed text" cell style, instead of
# minting a brand new number-format style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 1)
$ws.Range("I1").Value = "RMSInjectedVoltage"
$ws.Range("J1").Value = "RMSInjectedCurrent"

# New data row (row 2) - stored as text, matching sibling columns G2/H2
$ws.Range("I2").Value = "'60"
$ws.Range("J2").Value = "'1"

# Resize the new columns (and re-fit the now-adjacent column H) to their
# contents, mirroring the existing bestFit columns A:G.
$ws.Columns.Item(8).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(9).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(10).EntireColumn.AutoFit() | Out-Null

# Move the active selection to the next free cell, one column past the new
# data, matching where Excel would leave the cursor after data entry.
$ws.Range("K2").Select() | Out-Null
